$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revise constraints on basic asset classes portfolio (Max Weight column C)
$ws.Range("C3").Value = 0.75   # AGG
$ws.Range("C4").Value = 0.75   # IVV
$ws.Range("C6").Value = 0.3    # EFA
$ws.Range("C7").Value = 0.15   # EEM

# Update the active selection as recorded in the sheet view
$ws.Range("L17").Select()
